# Mark two wishlist items as reserved/bought ("Y" in column E):
#   - row 24: "Yoshi's Crafted World"
#   - row 26: "Set de potiuni"
# and leave the final selection on the last edited cell (E26), matching
# the user's click-through while ticking off items.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E24").Value = "Y"
$ws.Range("E26").Value = "Y"

$ws.Range("E26").Select()
